# Update the 100 arithmetic expressions in the "within 100" worksheet table.
# Each cell contains a simple "a+b=" / "a-b=" expression; replace the old
# expression text with the new one using Find/Replace on the whole document.
#
# wdFindContinue   = 1  (restart search from the top of the range each time)
# wdReplaceOne     = 1  (replace just the first/only match)
$d = $word.ActiveDocument

$replacements = @(
    @("64+34=", "51+16="),
    @("64-7=", "34+33="),
    @("38+1=", "73+21="),
    @("3-0=", "8+56="),
    @("81-42=", "7+54="),
    @("4+17=", "95-84="),
    @("74+24=", "62+27="),
    @("14+8=", "27+8="),
    @("71-51=", "8+12="),
    @("2+26=", "69+7="),
    @("23+70=", "45+45="),
    @("87-31=", "35-16="),
    @("43-8=", "22+55="),
    @("45+17=", "22+44="),
    @("46+3=", "85-63="),
    @("80-66=", "93-41="),
    @("19+26=", "76-32="),
    @("67+16=", "65+27="),
    @("50-45=", "69-14="),
    @("14+13=", "7+37="),
    @("35-33=", "79-0="),
    @("25+49=", "27+63="),
    @("70-40=", "66+2="),
    @("1+15=", "17+2="),
    @("29-21=", "30+38="),
    @("62-3=", "54+12="),
    @("78-77=", "71+22="),
    @("45+34=", "24-5="),
    @("16+31=", "54+17="),
    @("27+19=", "24+30="),
    @("62-24=", "43+20="),
    @("42+32=", "87-4="),
    @("22+35=", "25+28="),
    @("14-4=", "52-44="),
    @("1+69=", "98-92="),
    @("57+11=", "16+16="),
    @("70-63=", "29+43="),
    @("5+6=", "13-9="),
    @("64-8=", "15+6="),
    @("50+5=", "24-23="),
    @("27+60=", "98-65="),
    @("43-22=", "10+71="),
    @("76+15=", "48-13="),
    @("87-67=", "42+43="),
    @("21-20=", "2+97="),
    @("43+48=", "9+73="),
    @("44+13=", "65-38="),
    @("93-8=", "17+29="),
    @("74+13=", "1+5="),
    @("61+34=", "76-24="),
    @("37-14=", "66-3="),
    @("0+34=", "96-49="),
    @("61+13=", "1+88="),
    @("17+23=", "93-45="),
    @("10-7=", "52-48="),
    @("37-5=", "81-20="),
    @("31-16=", "28+27="),
    @("11+42=", "2+69="),
    @("27+53=", "4+4="),
    @("2+9=", "83-10="),
    @("60-26=", "34+29="),
    @("10+49=", "93-31="),
    @("78-74=", "32-32="),
    @("82-24=", "51+4="),
    @("84+8=", "19-12="),
    @("64-3=", "30+55="),
    @("64-40=", "91-65="),
    @("76+18=", "16+44="),
    @("45-5=", "72+11="),
    @("31-13=", "85-56="),
    @("90-63=", "51-16="),
    @("27+4=", "97-74="),
    @("50-6=", "0+17="),
    @("32+12=", "32+23="),
    @("11+8=", "16+52="),
    @("92-70=", "40+1="),
    @("40-10=", "89-87="),
    @("46-29=", "81-46="),
    @("52+31=", "18+12="),
    @("54+25=", "52-32="),
    @("76-25=", "34+28="),
    @("57-25=", "5+67="),
    @("65-31=", "55-36="),
    @("67-11=", "22+28="),
    @("2+42=", "9-1="),
    @("67-62=", "85-84="),
    @("61+28=", "46+14="),
    @("67-1=", "57-11="),
    @("37+48=", "11+48="),
    @("0+2=", "84-75="),
    @("84-28=", "75+24="),
    @("19+38=", "71-10="),
    @("85+8=", "23+69="),
    @("59-58=", "41+12="),
    @("32-4=", "76-63="),
    @("98-45=", "75+2="),
    @("4+31=", "42+7="),
    @("70-25=", "53+12="),
    @("68-9=", "82-2="),
    @("54-3=", "17+10=")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}
